$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4056.2856
$ws.Range("I74").Value = 3132.125
$ws.Range("K74").Value = 3132.125
$ws.Range("M74").Value = -2196.125
$ws.Range("H77").Value = 4056.2856
$ws.Range("I77").Value = 3132.125
$ws.Range("K77").Value = 15660.625
$ws.Range("M77").Value = -10980.625
$ws.Range("H135").Value = 992.2941
$ws.Range("I135").Value = 713.9091
$ws.Range("K135").Value = 6425.1819
$ws.Range("M135").Value = -3890.1819
$ws.Range("H138").Value = 9263619
$ws.Range("I138").Value = 1157.1875
$ws.Range("J138").Value = 16673589
$ws.Range("K138").Value = 3471.5625
$ws.Range("L138").Value = 50020767
$ws.Range("M138").Value = 1668.4375
$ws.Range("N138").Value = -50031047

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 43
$ws.Range("J4").Value = 30
$ws.Range("L4").Value = 30
$ws.Range("N4").Value = -262
$ws.Range("H61").Value = 4917.125
$ws.Range("I61").Value = 500
$ws.Range("K61").Value = 500
$ws.Range("M61").Value = -288
$ws.Range("H74").Value = 142667.75
$ws.Range("I74").Value = 188304.83
$ws.Range("K74").Value = 188304.83
$ws.Range("M74").Value = -187430.83
$ws.Range("H77").Value = 142667.75
$ws.Range("I77").Value = 188304.83
$ws.Range("K77").Value = 941524.1499999999
$ws.Range("M77").Value = -937156.1499999999
$ws.Range("H110").Value = 25589.533
$ws.Range("I110").Value = 28834.846
$ws.Range("J110").Value = 4495
$ws.Range("K110").Value = 28834.846
$ws.Range("L110").Value = 4495
$ws.Range("M110").Value = -26789.846
$ws.Range("N110").Value = -8585
$ws.Range("H132").Value = 4753
$ws.Range("I132").Value = 4720.375
$ws.Range("K132").Value = 14161.125
$ws.Range("M132").Value = -11631.125
$ws.Range("H134").Value = 45000
$ws.Range("J134").Value = 45000
$ws.Range("L134").Value = 45000
$ws.Range("N134").Value = -55140
$ws.Range("H135").Value = 69250
$ws.Range("J135").Value = 69250
$ws.Range("L135").Value = 69250
$ws.Range("N135").Value = -79390
$ws.Range("H136").Value = 4917.125
$ws.Range("I136").Value = 500
$ws.Range("K136").Value = 1500
$ws.Range("M136").Value = 1050
$ws.Range("H139").Value = 150715
$ws.Range("J139").Value = 150715
$ws.Range("L139").Value = 150715
$ws.Range("N139").Value = -160995
$ws.Range("H140").Value = 49999
$ws.Range("J140").Value = 49999
$ws.Range("L140").Value = 49999
$ws.Range("N140").Value = -60359

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 198.5
$ws.Range("I22").Value = 111.5
$ws.Range("J22").Value = 285.5
$ws.Range("K22").Value = 111.5
$ws.Range("L22").Value = 285.5
$ws.Range("M22").Value = 61.5
$ws.Range("N22").Value = -631.5
$ws.Range("H134").Value = 9099.4
$ws.Range("I134").Value = 8332.666999999999
$ws.Range("K134").Value = 24998.001
$ws.Range("M134").Value = -22463.001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 573.375
$ws.Range("I16").Value = 571.0769
$ws.Range("J16").Value = 583.3333
$ws.Range("K16").Value = 571.0769
$ws.Range("L16").Value = 583.3333
$ws.Range("M16").Value = -284.0769
$ws.Range("N16").Value = -1157.3333
$ws.Range("H105").Value = 16952.875
$ws.Range("I105").Value = 2583.3333
$ws.Range("J105").Value = 25574.6
$ws.Range("K105").Value = 2583.3333
$ws.Range("L105").Value = 25574.6
$ws.Range("M105").Value = -836.3332999999998
$ws.Range("N105").Value = -29068.6
$ws.Range("H113").Value = 573.375
$ws.Range("I113").Value = 571.0769
$ws.Range("J113").Value = 583.3333
$ws.Range("K113").Value = 571.0769
$ws.Range("L113").Value = 583.3333
$ws.Range("M113").Value = 1598.9231
$ws.Range("N113").Value = -4923.3333

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 800
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("H64").Value = 17942.715
$ws.Range("J64").Value = 17942.715
$ws.Range("L64").Value = 53828.145
$ws.Range("N64").Value = -54368.145
$ws.Range("H67").Value = 17942.715
$ws.Range("J67").Value = 17942.715
$ws.Range("L67").Value = 53828.145
$ws.Range("N67").Value = -55700.145
$ws.Range("H95").Value = 8600
$ws.Range("J95").Value = 8600
$ws.Range("L95").Value = 25800
$ws.Range("N95").Value = -29918
$ws.Range("H131").Value = 38035.484
$ws.Range("J131").Value = 5342.4
$ws.Range("L131").Value = 16027.2
$ws.Range("N131").Value = -26107.2
$ws.Range("H25").Value = 1009
$ws.Range("J25").Value = 1009
$ws.Range("L25").Value = 1009
$ws.Range("N25").Value = -2067
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("L35").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4148.364
$ws.Range("I132").Value = 4527.2
$ws.Range("K132").Value = 13581.6
$ws.Range("M132").Value = -11051.6

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4009.3125
$ws.Range("I136").Value = 1488.4117
$ws.Range("K136").Value = 4465.2351
$ws.Range("M136").Value = -1915.2351
$ws.Range("H137").Value = 75350.836
$ws.Range("J137").Value = 85428.75
$ws.Range("L137").Value = 85428.75
$ws.Range("N137").Value = -95628.75
$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360
$ws.Range("H141").Value = 73999.5
$ws.Range("J141").Value = 73999.5
$ws.Range("L141").Value = 73999.5
$ws.Range("N141").Value = -84359.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 675.875
$ws.Range("J81").Value = 824.5
$ws.Range("L81").Value = 1649
$ws.Range("N81").Value = -3771
$ws.Range("H84").Value = 675.875
$ws.Range("J84").Value = 824.5
$ws.Range("L84").Value = 8245
$ws.Range("N84").Value = -18853
$ws.Range("H107").Value = 292.07693
$ws.Range("I107").Value = 245.27272
$ws.Range("J107").Value = 549.5
$ws.Range("K107").Value = 735.81816
$ws.Range("L107").Value = 1648.5
$ws.Range("M107").Value = 1184.18184
$ws.Range("N107").Value = -5488.5
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("H126").Value = 6323.5386
$ws.Range("I126").Value = 6323.5386
$ws.Range("K126").Value = 18970.6158
$ws.Range("M126").Value = -16500.6158
$ws.Range("H136").Value = 1641.742
$ws.Range("I136").Value = 1446.7894
$ws.Range("K136").Value = 4340.3682
$ws.Range("M136").Value = -1790.3682
$ws.Range("H137").Value = 134998.67
$ws.Range("J137").Value = 134998.67
$ws.Range("L137").Value = 134998.67
$ws.Range("N137").Value = -145198.67
$ws.Range("H138").Value = 49999
$ws.Range("J138").Value = 49999
$ws.Range("L138").Value = 49999
$ws.Range("N138").Value = -60279
$ws.Range("H139").Value = 44999.5
$ws.Range("I139").Value = 40000
$ws.Range("J139").Value = 49999
$ws.Range("K139").Value = 40000
$ws.Range("L139").Value = 49999
$ws.Range("M139").Value = -34860
$ws.Range("N139").Value = -60279
$ws.Range("H141").Value = 98787.8
$ws.Range("J141").Value = 129999.5
$ws.Range("L141").Value = 129999.5
$ws.Range("N141").Value = -140359.5
$ws.Range("M108").ClearContents()
